$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data taken at mid mount: overwrite old tau values with -1 placeholders
$ws.Range("A1:E1").Value = -1
